# Weekly fruit/vegetable price update.
# A new daily record is inserted as row 33 (pushing the existing rows 33-101
# down to 34-102), matching the canonical diff: dimension grows from
# A1:R101 to A1:R102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; everything below (old rows 33-101)
# shifts down to 34-102 automatically, carrying its formatting (e.g. the
# date style on column D) along with it.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row with the new reading.
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44498
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112003
$ws.Range("G33").Value = "Ajo"
$ws.Range("H33").Value = "Chino"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 350
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 18000
$ws.Range("M33").Value = 17571
$ws.Range("N33").Value = "$/caja 10 kilos"
$ws.Range("O33").Value = "China"
$ws.Range("P33").Value = 1757
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"
